$wb = $excel.ActiveWorkbook

# Sheet1 row 6
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(6, 8).Value = 574.19354
$ws.Cells.Item(6, 9).Value = 107.69231
$ws.Cells.Item(6, 10).Value = 3000
$ws.Cells.Item(6, 11).Value = 323.07693
$ws.Cells.Item(6, 12).Value = 9000
$ws.Cells.Item(6, 13).Value = -211.07693
$ws.Cells.Item(6, 14).Value = -9224

# Sheet1 row 15
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(15, 8).Value = 938.02563
$ws.Cells.Item(15, 9).Value = 938.02563
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 2814.07689
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = -2645.07689

# Sheet1 row 29
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(29, 8).Value = 1500
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 1500
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 12).Value = 4500
$ws.Cells.Item(29, 13).ClearContents()
$ws.Cells.Item(29, 14).Value = -5062

# Sheet1 row 31
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(31, 8).Value = 2281.5
$ws.Cells.Item(31, 9).Value = 589
$ws.Cells.Item(31, 10).Value = 2620
$ws.Cells.Item(31, 11).Value = 1767
$ws.Cells.Item(31, 12).Value = 7860
$ws.Cells.Item(31, 13).Value = -1537
$ws.Cells.Item(31, 14).Value = -8320

# Sheet1 row 43
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(43, 8).Value = 1697.4667
$ws.Cells.Item(43, 9).Value = 925
$ws.Cells.Item(43, 10).Value = 2580.2856
$ws.Cells.Item(43, 11).Value = 925
$ws.Cells.Item(43, 12).Value = 2580.2856
$ws.Cells.Item(43, 13).Value = -856
$ws.Cells.Item(43, 14).Value = -2718.2856

# Sheet1 row 100
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(100, 8).Value = 2373.158
$ws.Cells.Item(100, 9).Value = 2108.182
$ws.Cells.Item(100, 10).Value = 2737.5
$ws.Cells.Item(100, 11).Value = 2108.182
$ws.Cells.Item(100, 12).Value = 2737.5
$ws.Cells.Item(100, 13).Value = -1567.182
$ws.Cells.Item(100, 14).Value = -3819.5

# Sheet1 row 116
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(116, 8).Value = 2508.2
$ws.Cells.Item(116, 9).Value = 2389.1667
$ws.Cells.Item(116, 10).Value = 2814.2856
$ws.Cells.Item(116, 11).Value = 2389.1667
$ws.Cells.Item(116, 12).Value = 2814.2856
$ws.Cells.Item(116, 13).Value = 1052.8333
$ws.Cells.Item(116, 14).Value = -9698.285599999999

# Sheet1 row 132
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(132, 8).Value = 1439.68
$ws.Cells.Item(132, 9).Value = 1044.0435
$ws.Cells.Item(132, 10).Value = 5989.5
$ws.Cells.Item(132, 11).Value = 3132.1305
$ws.Cells.Item(132, 12).Value = 17968.5
$ws.Cells.Item(132, 13).Value = -602.1305000000002
$ws.Cells.Item(132, 14).Value = -23028.5

# Sheet1 row 137
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(137, 8).Value = 3826.6316
$ws.Cells.Item(137, 9).Value = 4446.154
$ws.Cells.Item(137, 10).Value = 2484.3333
$ws.Cells.Item(137, 11).Value = 13338.462
$ws.Cells.Item(137, 12).Value = 7452.999899999999
$ws.Cells.Item(137, 13).Value = -10788.462
$ws.Cells.Item(137, 14).Value = -12552.9999

# Sheet1 row 138
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(138, 8).Value = 131898.88
$ws.Cells.Item(138, 9).Value = 1917.12
$ws.Cells.Item(138, 10).Value = 186975.89
$ws.Cells.Item(138, 11).Value = 5751.36
$ws.Cells.Item(138, 12).Value = 560927.67
$ws.Cells.Item(138, 13).Value = -611.3599999999997
$ws.Cells.Item(138, 14).Value = -571207.67

# Sheet2 row 4
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(4, 8).Value = 299
$ws.Cells.Item(4, 9).Value = 299
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 299
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = -183

# Sheet2 row 37
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(37, 8).Value = 12854.1
$ws.Cells.Item(37, 9).Value = 0
$ws.Cells.Item(37, 10).Value = 12854.1
$ws.Cells.Item(37, 11).Value = 0
$ws.Cells.Item(37, 12).Value = 12854.1
$ws.Cells.Item(37, 14).Value = -13400.1

# Sheet2 row 74
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(74, 8).Value = 1447.9375
$ws.Cells.Item(74, 9).Value = 1089.7693
$ws.Cells.Item(74, 10).Value = 3000
$ws.Cells.Item(74, 11).Value = 1089.7693
$ws.Cells.Item(74, 12).Value = 3000
$ws.Cells.Item(74, 13).Value = -215.7692999999999
$ws.Cells.Item(74, 14).Value = -4748

# Sheet2 row 77
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(77, 8).Value = 1447.9375
$ws.Cells.Item(77, 9).Value = 1089.7693
$ws.Cells.Item(77, 10).Value = 3000
$ws.Cells.Item(77, 11).Value = 5448.8465
$ws.Cells.Item(77, 12).Value = 15000
$ws.Cells.Item(77, 13).Value = -1080.8465
$ws.Cells.Item(77, 14).Value = -23736

# Sheet2 row 80
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(80, 8).Value = 21249
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 21249
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 21249
$ws.Cells.Item(80, 14).Value = -23245

# Sheet2 row 83
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(83, 8).Value = 21249
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 21249
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 63747
$ws.Cells.Item(83, 14).Value = -73731

# Sheet2 row 123
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(123, 8).Value = 40429
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 40429
$ws.Cells.Item(123, 11).Value = 0
$ws.Cells.Item(123, 12).Value = 40429
$ws.Cells.Item(123, 14).Value = -50229

# Sheet3 row 22
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(22, 8).Value = 3597.0625
$ws.Cells.Item(22, 9).Value = 3810.2
$ws.Cells.Item(22, 10).Value = 400
$ws.Cells.Item(22, 11).Value = 3810.2
$ws.Cells.Item(22, 12).Value = 400
$ws.Cells.Item(22, 13).Value = -3637.2
$ws.Cells.Item(22, 14).Value = -746

# Sheet3 row 80
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(80, 8).Value = 1219.2778
$ws.Cells.Item(80, 9).Value = 2420.5
$ws.Cells.Item(80, 10).Value = 258.3
$ws.Cells.Item(80, 11).Value = 2420.5
$ws.Cells.Item(80, 12).Value = 258.3
$ws.Cells.Item(80, 13).Value = -1422.5
$ws.Cells.Item(80, 14).Value = -2254.3

# Sheet3 row 83
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(83, 8).Value = 1219.2778
$ws.Cells.Item(83, 9).Value = 2420.5
$ws.Cells.Item(83, 10).Value = 258.3
$ws.Cells.Item(83, 11).Value = 12102.5
$ws.Cells.Item(83, 12).Value = 1291.5
$ws.Cells.Item(83, 13).Value = -7110.5
$ws.Cells.Item(83, 14).Value = -11275.5

# Sheet3 row 134
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(134, 8).Value = 3360.24
$ws.Cells.Item(134, 9).Value = 3050.3
$ws.Cells.Item(134, 10).Value = 4600
$ws.Cells.Item(134, 11).Value = 9150.900000000001
$ws.Cells.Item(134, 12).Value = 13800
$ws.Cells.Item(134, 13).Value = -6615.900000000001
$ws.Cells.Item(134, 14).Value = -18870

# Sheet4 row 31
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 4565.067
$ws.Cells.Item(31, 9).Value = 1347.75
$ws.Cells.Item(31, 10).Value = 5735
$ws.Cells.Item(31, 11).Value = 1347.75
$ws.Cells.Item(31, 12).Value = 5735
$ws.Cells.Item(31, 13).Value = -1052.75
$ws.Cells.Item(31, 14).Value = -6325

# Sheet4 row 34
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(34, 8).Value = 4565.067
$ws.Cells.Item(34, 9).Value = 1347.75
$ws.Cells.Item(34, 10).Value = 5735
$ws.Cells.Item(34, 11).Value = 1347.75
$ws.Cells.Item(34, 12).Value = 5735
$ws.Cells.Item(34, 13).Value = -1145.75
$ws.Cells.Item(34, 14).Value = -6139

# Sheet4 row 58
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(58, 8).Value = 1363.6875
$ws.Cells.Item(58, 9).Value = 1040.625
$ws.Cells.Item(58, 10).Value = 1686.75
$ws.Cells.Item(58, 11).Value = 1040.625
$ws.Cells.Item(58, 12).Value = 1686.75
$ws.Cells.Item(58, 13).Value = -837.625
$ws.Cells.Item(58, 14).Value = -2092.75

# Sheet4 row 105
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(105, 8).Value = 470
$ws.Cells.Item(105, 9).Value = 470
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 470
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).Value = 1277
$ws.Cells.Item(105, 14).ClearContents()

# Sheet4 row 136
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(136, 8).Value = 1363.6875
$ws.Cells.Item(136, 9).Value = 1040.625
$ws.Cells.Item(136, 10).Value = 1686.75
$ws.Cells.Item(136, 11).Value = 3121.875
$ws.Cells.Item(136, 12).Value = 5060.25
$ws.Cells.Item(136, 13).Value = -571.875
$ws.Cells.Item(136, 14).Value = -10160.25

# Sheet5 row 5
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(5, 8).Value = 1488.0294
$ws.Cells.Item(5, 9).Value = 516.7368
$ws.Cells.Item(5, 10).Value = 2718.3333
$ws.Cells.Item(5, 11).Value = 1550.2104
$ws.Cells.Item(5, 12).Value = 8154.999899999999
$ws.Cells.Item(5, 13).Value = -1438.2104
$ws.Cells.Item(5, 14).Value = -8378.999899999999

# Sheet5 row 122
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(122, 8).Value = 8974.5
$ws.Cells.Item(122, 9).Value = 588.3333
$ws.Cells.Item(122, 10).Value = 34133
$ws.Cells.Item(122, 11).Value = 5294.9997
$ws.Cells.Item(122, 12).Value = 307197
$ws.Cells.Item(122, 13).Value = -2844.9997
$ws.Cells.Item(122, 14).Value = -312097

# Sheet5 row 132
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(132, 8).Value = 3660.9666
$ws.Cells.Item(132, 9).Value = 2615.36
$ws.Cells.Item(132, 10).Value = 4407.8286
$ws.Cells.Item(132, 11).Value = 23538.24
$ws.Cells.Item(132, 12).Value = 39670.4574
$ws.Cells.Item(132, 13).Value = -21008.24
$ws.Cells.Item(132, 14).Value = -44730.4574

# Sheet5 row 135
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(135, 8).Value = 1488.0294
$ws.Cells.Item(135, 9).Value = 516.7368
$ws.Cells.Item(135, 10).Value = 2718.3333
$ws.Cells.Item(135, 11).Value = 4650.6312
$ws.Cells.Item(135, 12).Value = 24464.9997
$ws.Cells.Item(135, 13).Value = -2115.6312
$ws.Cells.Item(135, 14).Value = -29534.9997

# Sheet6 row 2
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 8).Value = 63.35294
$ws.Cells.Item(2, 9).Value = 64.933334
$ws.Cells.Item(2, 10).Value = 51.5
$ws.Cells.Item(2, 11).Value = 64.933334
$ws.Cells.Item(2, 12).Value = 51.5
$ws.Cells.Item(2, 13).Value = 48.066666
$ws.Cells.Item(2, 14).Value = -277.5

# Sheet6 row 46
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(46, 8).Value = 4499.9443
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 10).Value = 4499.9443
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 12).Value = 4499.9443
$ws.Cells.Item(46, 14).Value = -4811.9443

# Sheet6 row 80
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(80, 8).Value = 603372.25
$ws.Cells.Item(80, 9).Value = 1289597
$ws.Cells.Item(80, 10).Value = 2925.625
$ws.Cells.Item(80, 11).Value = 1289597
$ws.Cells.Item(80, 12).Value = 2925.625
$ws.Cells.Item(80, 13).Value = -1288599
$ws.Cells.Item(80, 14).Value = -4921.625

# Sheet6 row 83
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(83, 8).Value = 603372.25
$ws.Cells.Item(83, 9).Value = 1289597
$ws.Cells.Item(83, 10).Value = 2925.625
$ws.Cells.Item(83, 11).Value = 6447985
$ws.Cells.Item(83, 12).Value = 14628.125
$ws.Cells.Item(83, 13).Value = -6442993
$ws.Cells.Item(83, 14).Value = -24612.125

# Sheet6 row 122
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(122, 8).Value = 4484.387
$ws.Cells.Item(122, 9).Value = 2986.3333
$ws.Cells.Item(122, 10).Value = 5430.5264
$ws.Cells.Item(122, 11).Value = 8958.999899999999
$ws.Cells.Item(122, 12).Value = 16291.5792
$ws.Cells.Item(122, 13).Value = -6508.999899999999
$ws.Cells.Item(122, 14).Value = -21191.5792

# Sheet6 row 132
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(132, 8).Value = 2907.625
$ws.Cells.Item(132, 9).Value = 2252.4
$ws.Cells.Item(132, 10).Value = 3999.6667
$ws.Cells.Item(132, 11).Value = 6757.200000000001
$ws.Cells.Item(132, 12).Value = 11999.0001
$ws.Cells.Item(132, 13).Value = -4227.200000000001
$ws.Cells.Item(132, 14).Value = -17059.0001

# Sheet7 row 7
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(7, 8).Value = 43480732
$ws.Cells.Item(7, 9).Value = 55558144
$ws.Cells.Item(7, 10).Value = 2061
$ws.Cells.Item(7, 11).Value = 55558144
$ws.Cells.Item(7, 12).Value = 2061
$ws.Cells.Item(7, 13).Value = -55558032
$ws.Cells.Item(7, 14).Value = -2285

# Sheet7 row 74
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(74, 8).Value = 33000
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 33000
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 33000
$ws.Cells.Item(74, 13).ClearContents()
$ws.Cells.Item(74, 14).Value = -34996

# Sheet7 row 77
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(77, 8).Value = 33000
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 33000
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 12).Value = 99000
$ws.Cells.Item(77, 13).ClearContents()
$ws.Cells.Item(77, 14).Value = -108984

# Sheet7 row 126
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(126, 8).Value = 43480732
$ws.Cells.Item(126, 9).Value = 55558144
$ws.Cells.Item(126, 10).Value = 2061
$ws.Cells.Item(126, 11).Value = 166674432
$ws.Cells.Item(126, 12).Value = 6183
$ws.Cells.Item(126, 13).Value = -166671962
$ws.Cells.Item(126, 14).Value = -11123

# Sheet7 row 132
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(132, 8).Value = 3784.182
$ws.Cells.Item(132, 9).Value = 3270.3333
$ws.Cells.Item(132, 10).Value = 5154.4443
$ws.Cells.Item(132, 11).Value = 9810.999899999999
$ws.Cells.Item(132, 12).Value = 15463.3329
$ws.Cells.Item(132, 13).Value = -7280.999899999999
$ws.Cells.Item(132, 14).Value = -20523.3329

# Sheet7 row 136
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(136, 8).Value = 11112973
$ws.Cells.Item(136, 9).Value = 2212.125
$ws.Cells.Item(136, 10).Value = 23810986
$ws.Cells.Item(136, 11).Value = 6636.375
$ws.Cells.Item(136, 12).Value = 71432958
$ws.Cells.Item(136, 13).Value = -4086.375
$ws.Cells.Item(136, 14).Value = -71438058

# Sheet8 row 98
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(98, 8).Value = 90000
$ws.Cells.Item(98, 9).Value = 0
$ws.Cells.Item(98, 10).Value = 90000
$ws.Cells.Item(98, 11).Value = 0
$ws.Cells.Item(98, 12).Value = 90000
$ws.Cells.Item(98, 14).Value = -95990

# Sheet8 row 118
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(118, 8).Value = 50392.5
$ws.Cells.Item(118, 9).Value = 0
$ws.Cells.Item(118, 10).Value = 50392.5
$ws.Cells.Item(118, 11).Value = 0
$ws.Cells.Item(118, 12).Value = 50392.5
$ws.Cells.Item(118, 14).Value = -53706.5

# Sheet8 row 123
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(123, 8).Value = 34194.77
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 34194.77
$ws.Cells.Item(123, 11).Value = 0
$ws.Cells.Item(123, 12).Value = 34194.77
$ws.Cells.Item(123, 14).Value = -43994.77

# Sheet8 row 136
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(136, 8).Value = 3197.6843
$ws.Cells.Item(136, 9).Value = 2909.0417
$ws.Cells.Item(136, 10).Value = 3692.5
$ws.Cells.Item(136, 11).Value = 8727.125100000001
$ws.Cells.Item(136, 12).Value = 11077.5
$ws.Cells.Item(136, 13).Value = -6177.125100000001
$ws.Cells.Item(136, 14).Value = -16177.5
